# edit.ps1 - apply weekly CompStat crime-data refresh (NYPD 014 Pct week of 6/26/2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/report-week banner) ---
$ws.Range("A8").Value = "Volume 30   Number  26"
$ws.Range("C9").Value = "Report Covering the Week  6/26/2023  Through  7/2/2023"

# --- Straightforward numeric value updates (type/style unchanged) ---
$ws.Range("N14").Value = -71.428571428571
$ws.Range("L15").Value = -22.222222222222
$ws.Range("M15").Value = 133.333333333333
$ws.Range("N15").Value = -58.823529411764
$ws.Range("F16").Value = 39
$ws.Range("G16").Value = 58
$ws.Range("H16").Value = -32.758620689655
$ws.Range("I16").Value = 246
$ws.Range("J16").Value = 298
$ws.Range("K16").Value = -17.44966442953
$ws.Range("L16").Value = 36.666666666666
$ws.Range("M16").Value = 267.164179104478
$ws.Range("N16").Value = -80.811232449298
$ws.Range("C17").Value = 12
$ws.Range("D17").Value = 14
$ws.Range("E17").Value = -14.285714285714
$ws.Range("F17").Value = 29
$ws.Range("G17").Value = 50
$ws.Range("H17").Value = -42
$ws.Range("I17").Value = 239
$ws.Range("J17").Value = 207
$ws.Range("K17").Value = 15.458937198067
$ws.Range("L17").Value = 4.366812227074
$ws.Range("M17").Value = 151.578947368421
$ws.Range("N17").Value = -25.776397515528
$ws.Range("C18").Value = 9
$ws.Range("D18").Value = 20
$ws.Range("E18").Value = -55
$ws.Range("F18").Value = 24
$ws.Range("G18").Value = 71
$ws.Range("H18").Value = -66.197183098591
$ws.Range("I18").Value = 219
$ws.Range("J18").Value = 356
$ws.Range("K18").Value = -38.483146067415
$ws.Range("L18").Value = 9.5
$ws.Range("M18").Value = 26.589595375722
$ws.Range("N18").Value = -83.496608892238
$ws.Range("C19").Value = 42
$ws.Range("D19").Value = 61
$ws.Range("E19").Value = -31.147540983606
$ws.Range("F19").Value = 185
$ws.Range("G19").Value = 212
$ws.Range("H19").Value = -12.735849056603
$ws.Range("I19").Value = 1218
$ws.Range("J19").Value = 1072
$ws.Range("K19").Value = 13.619402985074
$ws.Range("L19").Value = 111.091854419411
$ws.Range("M19").Value = 7.029876977152
$ws.Range("N19").Value = -74.05750798722
$ws.Range("C20").Value = 2
$ws.Range("E20").Value = 100
$ws.Range("I20").Value = 28
$ws.Range("J20").Value = 40
$ws.Range("K20").Value = -30
$ws.Range("L20").Value = 16.666666666666
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = -84.782608695652
$ws.Range("D21").Value = 114
$ws.Range("E21").Value = -35.964912280701
$ws.Range("F21").Value = 281
$ws.Range("G21").Value = 400
$ws.Range("H21").Value = -29.75
$ws.Range("I21").Value = 1959
$ws.Range("J21").Value = 1989
$ws.Range("K21").Value = -1.508295625942
$ws.Range("L21").Value = 60.44226044226
$ws.Range("M21").Value = 31.476510067114
$ws.Range("N21").Value = -74.993617564462
$ws.Range("C22").Value = 3
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = -25
$ws.Range("F22").Value = 16
$ws.Range("G22").Value = 13
$ws.Range("H22").Value = 23.076923076923
$ws.Range("I22").Value = 114
$ws.Range("J22").Value = 96
$ws.Range("K22").Value = 18.75
$ws.Range("L22").Value = 67.647058823529
$ws.Range("M22").Value = 72.727272727272
$ws.Range("C24").Value = 76
$ws.Range("D24").Value = 73
$ws.Range("E24").Value = 4.109589041095
$ws.Range("F24").Value = 379
$ws.Range("G24").Value = 278
$ws.Range("H24").Value = 36.330935251798
$ws.Range("I24").Value = 1995
$ws.Range("J24").Value = 1524
$ws.Range("K24").Value = 30.905511811023
$ws.Range("L24").Value = 92.940038684719
$ws.Range("M24").Value = -15.322580645161
$ws.Range("C25").Value = 24
$ws.Range("D25").Value = 25
$ws.Range("E25").Value = -4
$ws.Range("F25").Value = 93
$ws.Range("G25").Value = 91
$ws.Range("H25").Value = 2.197802197802
$ws.Range("I25").Value = 514
$ws.Range("J25").Value = 452
$ws.Range("K25").Value = 13.716814159292
$ws.Range("L25").Value = 18.160919540229
$ws.Range("M25").Value = 76.632302405498
$ws.Range("F26").Value = 1
$ws.Range("L26").Value = -35
$ws.Range("C27").Value = 6
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 200
$ws.Range("F27").Value = 26
$ws.Range("G27").Value = 23
$ws.Range("H27").Value = 13.043478260869
$ws.Range("I27").Value = 116
$ws.Range("J27").Value = 106
$ws.Range("K27").Value = 9.43396226415
$ws.Range("L27").Value = 70.588235294117
$ws.Range("J28").Value = 4
$ws.Range("K28").Value = -25
$ws.Range("L28").Value = -66.666666666666
$ws.Range("J29").Value = 2
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = -71.428571428571
$ws.Range("I30").Value = 7
$ws.Range("K30").Value = -41.666666666666
$ws.Range("L30").Value = -65

# --- Cells changing from blank-marker text ("0"/"***.*") to real numeric values ---
# Use Copy + two-stage PasteSpecial (Formats, then Values) from a same-style numeric
# donor cell so the destination picks up the correct style index, then overwrite the
# pasted value with the real figure.
function Set-NumericFromText($destAddr, $donorAddr, $value) {
    $ws.Range($donorAddr).Copy()
    $ws.Range($destAddr).PasteSpecial(-4122)
    $ws.Range($donorAddr).Copy()
    $ws.Range($destAddr).PasteSpecial(-4163)
    $excel.CutCopyMode = $false
    $ws.Range($destAddr).Value = $value
}

# style15 (#,##0) donor: I28 (stays 3, untouched by this edit)
Set-NumericFromText "D28" "I28" 2
Set-NumericFromText "G28" "I28" 2
Set-NumericFromText "D29" "I28" 1
Set-NumericFromText "G29" "I28" 1
Set-NumericFromText "F30" "I28" 1

# style16 (#,##0.0 pct) donor: K14 (stays -60, untouched by this edit)
Set-NumericFromText "E28" "K14" -100
Set-NumericFromText "H28" "K14" -100
Set-NumericFromText "E29" "K14" -100
Set-NumericFromText "H29" "K14" -100

# --- Cells changing from a real numeric value back to the blank-marker text ---
# Donor cells D15 ("0") / E15 ("***.*") keep shared-string text + style14 intact,
# and are themselves untouched by this edit, so they are safe, stable donors.
function Set-TextFromDonor($destAddr, $donorAddr) {
    $ws.Range($donorAddr).Copy()
    $ws.Range($destAddr).PasteSpecial(-4122)
    $ws.Range($donorAddr).Copy()
    $ws.Range($destAddr).PasteSpecial(-4163)
    $excel.CutCopyMode = $false
}

Set-TextFromDonor "C15" "D15"
Set-TextFromDonor "C26" "D15"
Set-TextFromDonor "G30" "D15"
Set-TextFromDonor "H30" "E15"
